$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.929.04"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.306.70"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("D9").Value = "2.304.52"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "59.926.03"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "2.718.33"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D18").Value = "2.320.73"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "0.0₃0725"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.378"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "316.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.569"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("E49").Value = "  +21.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0213"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("E51").Value = "  +0.19%  "
